$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 607663.9
$ws.Range("I132").Value = 2510.8872
$ws.Range("J132").Value = 4904250
$ws.Range("K132").Value = 7532.6616
$ws.Range("L132").Value = 14712750
$ws.Range("M132").Value = -5002.6616
$ws.Range("N132").Value = -14717810

# Row 137
$ws.Range("H137").Value = 2383111.5
$ws.Range("I137").Value = 4349706
$ws.Range("J137").Value = 2497.0527
$ws.Range("K137").Value = 13049118
$ws.Range("L137").Value = 7491.158100000001
$ws.Range("M137").Value = -13046568
$ws.Range("N137").Value = -12591.1581

# Row 138
$ws.Range("H138").Value = 3034086.8
$ws.Range("I138").Value = 2446.6191
$ws.Range("J138").Value = 4448852
$ws.Range("K138").Value = 7339.8573
$ws.Range("L138").Value = 13346556
$ws.Range("M138").Value = -2199.8573
$ws.Range("N138").Value = -13356836


$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 43567130
$ws.Range("I61").Value = 52686176
$ws.Range("J61").Value = 251650
$ws.Range("K61").Value = 52686176
$ws.Range("L61").Value = 251650
$ws.Range("M61").Value = -52685964
$ws.Range("N61").Value = -252074

# Row 131
$ws.Range("H131").Value = 55041.668
$ws.Range("J131").Value = 55041.668
$ws.Range("L131").Value = 55041.668
$ws.Range("N131").Value = -65121.668

# Row 136
$ws.Range("H136").Value = 43567130
$ws.Range("I136").Value = 52686176
$ws.Range("J136").Value = 251650
$ws.Range("K136").Value = 158058528
$ws.Range("L136").Value = 754950
$ws.Range("M136").Value = -158055978
$ws.Range("N136").Value = -760050


$ws = $wb.Worksheets.Item("BSM")
# Row 130
$ws.Range("H130").Value = 46102
$ws.Range("J130").Value = 46102
$ws.Range("L130").Value = 46102
$ws.Range("N130").Value = -56142

# Row 134
$ws.Range("H134").Value = 13336835
$ws.Range("I134").Value = 3489.375
$ws.Range("J134").Value = 37040560
$ws.Range("K134").Value = 10468.125
$ws.Range("L134").Value = 111121680
$ws.Range("M134").Value = -7933.125
$ws.Range("N134").Value = -111126750


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 56860.95
$ws.Range("I31").Value = 34410.902
$ws.Range("J31").Value = 134188.89
$ws.Range("K31").Value = 34410.902
$ws.Range("L31").Value = 134188.89
$ws.Range("M31").Value = -34115.902
$ws.Range("N31").Value = -134778.89

# Row 34
$ws.Range("H34").Value = 56860.95
$ws.Range("I34").Value = 34410.902
$ws.Range("J34").Value = 134188.89
$ws.Range("K34").Value = 34410.902
$ws.Range("L34").Value = 134188.89
$ws.Range("M34").Value = -34208.902
$ws.Range("N34").Value = -134592.89

# Row 130
$ws.Range("H130").Value = 57500
$ws.Range("J130").Value = 57500
$ws.Range("L130").Value = 57500
$ws.Range("N130").Value = -67540

# Row 132
$ws.Range("H132").Value = 32365
$ws.Range("I132").Value = 1861.2273
$ws.Range("J132").Value = 93372.55
$ws.Range("K132").Value = 5583.6819
$ws.Range("L132").Value = 280117.65
$ws.Range("M132").Value = -3053.6819
$ws.Range("N132").Value = -285177.65

# Row 134
$ws.Range("H134").Value = 35856.03
$ws.Range("I134").Value = 911.43475
$ws.Range("J134").Value = 136321.75
$ws.Range("K134").Value = 2734.30425
$ws.Range("L134").Value = 408965.25
$ws.Range("M134").Value = -199.3042500000001
$ws.Range("N134").Value = -414035.25


$ws = $wb.Worksheets.Item("CUL")
# Row 130
$ws.Range("H130").Value = 2852.1
$ws.Range("I130").Value = 965.5
$ws.Range("K130").Value = 2896.5
$ws.Range("M130").Value = 2123.5

# Row 131
$ws.Range("H131").Value = 1037.5
$ws.Range("J131").Value = 1079.2683
$ws.Range("L131").Value = 3237.8049
$ws.Range("N131").Value = -13317.8049


$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2659.9333
$ws.Range("I102").Value = 2707.0833
$ws.Range("K102").Value = 2707.0833
$ws.Range("M102").Value = -1085.0833

# Row 113
$ws.Range("H113").Value = 2457.75
$ws.Range("I113").Value = 1392.6666
$ws.Range("J113").Value = 3827.1428
$ws.Range("K113").Value = 1392.6666
$ws.Range("L113").Value = 3827.1428
$ws.Range("M113").Value = 777.3334
$ws.Range("N113").Value = -8167.1428

# Row 132
$ws.Range("H132").Value = 74071.18
$ws.Range("I132").Value = 65046.312
$ws.Range("J132").Value = 86104.336
$ws.Range("K132").Value = 195138.936
$ws.Range("L132").Value = 258313.008
$ws.Range("M132").Value = -192608.936
$ws.Range("N132").Value = -263373.008


$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2700
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2364
$ws.Range("N40").Value = -3272

# Row 61
$ws.Range("H61").Value = 2580.08
$ws.Range("I61").Value = 2594.3157
$ws.Range("J61").Value = 2535
$ws.Range("K61").Value = 2594.3157
$ws.Range("L61").Value = 2535
$ws.Range("M61").Value = -2392.3157
$ws.Range("N61").Value = -2939

# Row 82
$ws.Range("H82").Value = 1316.2307
$ws.Range("I82").Value = 1097.4286
$ws.Range("J82").Value = 1571.5
$ws.Range("K82").Value = 1097.4286
$ws.Range("L82").Value = 1571.5
$ws.Range("M82").Value = -736.4286
$ws.Range("N82").Value = -2293.5

# Row 85
$ws.Range("H85").Value = 1316.2307
$ws.Range("I85").Value = 1097.4286
$ws.Range("J85").Value = 1571.5
$ws.Range("K85").Value = 1097.4286
$ws.Range("L85").Value = 1571.5
$ws.Range("M85").Value = 150.5714
$ws.Range("N85").Value = -4067.5

# Row 93
$ws.Range("H93").Value = 1473.6666
$ws.Range("I93").Value = 1237
$ws.Range("J93").Value = 1947
$ws.Range("K93").Value = 1237
$ws.Range("L93").Value = 1947
$ws.Range("M93").Value = 11
$ws.Range("N93").Value = -4443

# Row 113
$ws.Range("H113").Value = 2580.08
$ws.Range("I113").Value = 2594.3157
$ws.Range("J113").Value = 2535
$ws.Range("K113").Value = 2594.3157
$ws.Range("L113").Value = 2535
$ws.Range("M113").Value = -424.3157000000001
$ws.Range("N113").Value = -6875

# Row 127
$ws.Range("H127").Value = 52500
$ws.Range("J127").Value = 52500
$ws.Range("L127").Value = 52500
$ws.Range("N127").Value = -62420

# Row 132
$ws.Range("H132").Value = 103384.8
$ws.Range("I132").Value = 2774.8333
$ws.Range("J132").Value = 254299.75
$ws.Range("K132").Value = 8324.499899999999
$ws.Range("L132").Value = 762899.25
$ws.Range("M132").Value = -5794.499899999999
$ws.Range("N132").Value = -767959.25


$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 723.7368
$ws.Range("I107").Value = 550.9231
$ws.Range("J107").Value = 1098.1666
$ws.Range("K107").Value = 1652.7693
$ws.Range("L107").Value = 3294.4998
$ws.Range("M107").Value = 267.2307000000001
$ws.Range("N107").Value = -7134.4998

# Row 122
$ws.Range("H122").Value = 2608.9
$ws.Range("I122").Value = 1983.1666
$ws.Range("J122").Value = 3547.5
$ws.Range("K122").Value = 5949.4998
$ws.Range("L122").Value = 10642.5
$ws.Range("M122").Value = -3499.4998
$ws.Range("N122").Value = -15542.5

# Row 125
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840

# Row 126
$ws.Range("H126").Value = 934.8
$ws.Range("I126").Value = 934.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2804.4
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -334.3999999999996
$ws.Range("N126").ClearContents()

